$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet "Test Cases" -> "TestCases" (this also updates the
# _xlnm._FilterDatabase defined name reference automatically).
$ws.Name = "TestCases"

# Update the zoom level of the sheet view from 55% to 70%.
$ws.Application.ActiveWindow.Zoom = 70

# Correzione TraceId test #6 e #147
# Test #6 is on row 8, TRACEID column is H.
$ws.Range("H8").Value = "d6ba896e12aecce9"

# Test #147 is on row 88, TRACEID column is H.
$ws.Range("H88").Value = "ae72ee4beb45bf2f"
